$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 234.74074
$ws.Range("I33").Value = 166.84616
$ws.Range("J33").Value = 2000
$ws.Range("K33").Value = 166.84616
$ws.Range("L33").Value = 2000
$ws.Range("M33").Value = 62.15384

$ws.Range("H40").Value = 2407908.2
$ws.Range("I40").Value = 6383.4
$ws.Range("J40").Value = 3908861
$ws.Range("K40").Value = 6383.4
$ws.Range("L40").Value = 3908861
$ws.Range("M40").Value = -6208.4
$ws.Range("N40").Value = -3909211

$ws.Range("H138").Value = 2014.6888
$ws.Range("I138").Value = 1472.25
$ws.Range("J138").Value = 2908.1177
$ws.Range("K138").Value = 4416.75
$ws.Range("L138").Value = 8724.3531
$ws.Range("M138").Value = 723.25
$ws.Range("N138").Value = -19004.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6063449.5
$ws.Range("I32").Value = 6063449.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 6063449.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -6063162.5

$ws.Range("H64").Value = 45000
$ws.Range("I64").Value = 45000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 45000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -44752
$ws.Range("N64").Value = $null

$ws.Range("H67").Value = 45000
$ws.Range("I67").Value = 45000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 45000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -44142
$ws.Range("N67").Value = $null

$ws.Range("H74").Value = 3637.16
$ws.Range("I74").Value = 1019.4706
$ws.Range("J74").Value = 9199.75
$ws.Range("K74").Value = 1019.4706
$ws.Range("L74").Value = 9199.75
$ws.Range("M74").Value = -145.4706

$ws.Range("H77").Value = 3637.16
$ws.Range("I77").Value = 1019.4706
$ws.Range("J77").Value = 9199.75
$ws.Range("K77").Value = 5097.353
$ws.Range("L77").Value = 45998.75
$ws.Range("M77").Value = -729.3530000000001

$ws.Range("H97").Value = 1925.5714
$ws.Range("I97").Value = 1746.5
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 1746.5
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1250.5
$ws.Range("N97").Value = -3992

$ws.Range("H110").Value = 1378.1538
$ws.Range("I110").Value = 1491.8
$ws.Range("J110").Value = 999.3333
$ws.Range("K110").Value = 1491.8
$ws.Range("L110").Value = 999.3333
$ws.Range("M110").Value = 553.2

$ws.Range("H132").Value = 819678.75
$ws.Range("I132").Value = 848810.1
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 2546430.3
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2543900.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 80000
$ws.Range("I62").Value = 80000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 80000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -79314
$ws.Range("N62").Value = $null

$ws.Range("H65").Value = 80000
$ws.Range("I65").Value = 80000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 240000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -236568
$ws.Range("N65").Value = $null

$ws.Range("H94").Value = 1424.8462
$ws.Range("I94").Value = 1132.7
$ws.Range("J94").Value = 2398.6667
$ws.Range("K94").Value = 1132.7
$ws.Range("L94").Value = 2398.6667
$ws.Range("M94").Value = -681.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 91812.02
$ws.Range("I31").Value = 157515.2
$ws.Range("J31").Value = 22650.79
$ws.Range("K31").Value = 157515.2
$ws.Range("L31").Value = 22650.79
$ws.Range("M31").Value = -157220.2
$ws.Range("N31").Value = -23240.79

$ws.Range("H34").Value = 91812.02
$ws.Range("I34").Value = 157515.2
$ws.Range("J34").Value = 22650.79
$ws.Range("K34").Value = 157515.2
$ws.Range("L34").Value = 22650.79
$ws.Range("M34").Value = -157313.2
$ws.Range("N34").Value = -23054.79

$ws.Range("H107").Value = 667.8461
$ws.Range("I107").Value = 516.5454999999999
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 516.5454999999999
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1403.4545
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 265.0625
$ws.Range("I2").Value = 512.8570999999999
$ws.Range("J2").Value = 72.333336
$ws.Range("K2").Value = 3077.1426
$ws.Range("L2").Value = 434.000016
$ws.Range("M2").Value = -2964.1426
$ws.Range("N2").Value = -660.000016

$ws.Range("H23").Value = 7812845.5
$ws.Range("I23").Value = 83
$ws.Range("J23").Value = 20834116
$ws.Range("K23").Value = 249
$ws.Range("L23").Value = 62502348
$ws.Range("M23").Value = -14
$ws.Range("N23").Value = -62502818

$ws.Range("H75").Value = 7050.4287
$ws.Range("I75").Value = 1463
$ws.Range("J75").Value = 9285.4
$ws.Range("K75").Value = 4389
$ws.Range("L75").Value = 27856.2
$ws.Range("M75").Value = -3391
$ws.Range("N75").Value = -29852.2

$ws.Range("H78").Value = 7050.4287
$ws.Range("I78").Value = 1463
$ws.Range("J78").Value = 9285.4
$ws.Range("K78").Value = 13167
$ws.Range("L78").Value = 83568.59999999999
$ws.Range("M78").Value = -8175
$ws.Range("N78").Value = -93552.59999999999

$ws.Range("H86").Value = 277.33334
$ws.Range("I86").Value = 274.25
$ws.Range("J86").Value = 279.8
$ws.Range("K86").Value = 822.75
$ws.Range("L86").Value = 839.4000000000001
$ws.Range("M86").Value = 363.25

$ws.Range("H89").Value = 277.33334
$ws.Range("I89").Value = 274.25
$ws.Range("J89").Value = 279.8
$ws.Range("K89").Value = 2468.25
$ws.Range("L89").Value = 2518.2
$ws.Range("M89").Value = 3459.75

$ws.Range("H125").Value = 19031.8

$ws.Range("H131").Value = 11710.479
$ws.Range("I131").Value = 1526.25
$ws.Range("J131").Value = 13854.526
$ws.Range("K131").Value = 4578.75
$ws.Range("L131").Value = 41563.578
$ws.Range("M131").Value = 461.25
$ws.Range("N131").Value = -51643.578

$ws.Range("H132").Value = 2480.818
$ws.Range("I132").Value = 1274.8334
$ws.Range("J132").Value = 3928
$ws.Range("K132").Value = 11473.5006
$ws.Range("L132").Value = 35352
$ws.Range("M132").Value = -8943.500599999999
$ws.Range("N132").Value = -40412

$ws.Range("H137").Value = 3638.1428
$ws.Range("I137").Value = 2128.4
$ws.Range("J137").Value = 7412.5
$ws.Range("K137").Value = 6385.200000000001
$ws.Range("L137").Value = 22237.5
$ws.Range("M137").Value = -1285.200000000001
$ws.Range("N137").Value = -32437.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 57925
$ws.Range("I39").Value = 16850
$ws.Range("J39").Value = 99000
$ws.Range("K39").Value = 16850
$ws.Range("L39").Value = 99000
$ws.Range("M39").Value = -16318
$ws.Range("N39").Value = -100064

$ws.Range("H58").Value = 31748.75
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 31748.75
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 31748.75
$ws.Range("N58").Value = -32302.75

$ws.Range("H70").Value = 6098.3076
$ws.Range("I70").Value = 6162.5454
$ws.Range("J70").Value = 5745
$ws.Range("K70").Value = 6162.5454
$ws.Range("L70").Value = 5745
$ws.Range("M70").Value = -5892.5454
$ws.Range("N70").Value = -6285

$ws.Range("H73").Value = 6098.3076
$ws.Range("I73").Value = 6162.5454
$ws.Range("J73").Value = 5745
$ws.Range("K73").Value = 6162.5454
$ws.Range("L73").Value = 5745
$ws.Range("M73").Value = -5226.5454
$ws.Range("N73").Value = -7617

$ws.Range("H80").Value = 2561561.5
$ws.Range("I80").Value = 2561561.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2561561.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2560563.5
$ws.Range("N80").Value = $null

$ws.Range("H83").Value = 2561561.5
$ws.Range("I83").Value = 2561561.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 12807807.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -12802815.5
$ws.Range("N83").Value = $null

$ws.Range("H102").Value = 2963.761
$ws.Range("I102").Value = 2433.575
$ws.Range("J102").Value = 6498.3335
$ws.Range("K102").Value = 2433.575
$ws.Range("L102").Value = 6498.3335
$ws.Range("M102").Value = -811.5749999999998
$ws.Range("N102").Value = -9742.333500000001

$ws.Range("H122").Value = 41267.43
$ws.Range("I122").Value = 65184.312
$ws.Range("J122").Value = 9378.25
$ws.Range("K122").Value = 195552.936
$ws.Range("L122").Value = 28134.75
$ws.Range("M122").Value = -193102.936
$ws.Range("N122").Value = -33034.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4885.25
$ws.Range("I40").Value = 4816.4
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 4816.4
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -4680.4

$ws.Range("H57").Value = 34999
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 34999
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 34999
$ws.Range("N57").Value = -36131

$ws.Range("H122").Value = 4987.9287
$ws.Range("I122").Value = 4893.9546
$ws.Range("J122").Value = 5332.5
$ws.Range("K122").Value = 14681.8638
$ws.Range("L122").Value = 15997.5
$ws.Range("M122").Value = -12231.8638
$ws.Range("N122").Value = -20897.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 29999
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 29999
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 29999
$ws.Range("N48").Value = -31137

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = $null

$ws.Range("H122").Value = 3183.8386
$ws.Range("I122").Value = 2499.7827
$ws.Range("J122").Value = 5150.5
$ws.Range("K122").Value = 7499.348100000001
$ws.Range("L122").Value = 15451.5
$ws.Range("M122").Value = -5049.348100000001
$ws.Range("N122").Value = -20351.5

